$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old "Terms Typically Offered" column (D),
# shifting it to G. This also shifts the sheet dimension from A1:D34 to A1:G34.
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# New header row labels for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default every data row's new Corequisites / Concurrent / Recommended cells to "NA".
$ws.Range("D2:F34").Value = "NA"

# Row 14 (WVIT 344): drop the "Co-requisite: " phrasing now that it has its own column.
$ws.Range("C14").Value = "WVIT or RPTA major. AGB 212 or ECON 201 or ECON 221 or ECON 222."

# Row 23 (WVIT 428): move the "Recommended:" clause out of Prerequisites and into the
# new Recommended column; the Terms cell also picks up a trailing space in the source.
$ws.Range("C23").Value = "AGB 214, SS 221 and WVIT 233."
$ws.Range("F23").Value = "WVIT 331, WVIT 332, and WVIT 333."
$ws.Range("G23").Value = "SP "

# Row 28 (WVIT 460): drop "the following" from the prerequisite text.
$ws.Range("C28").Value = "WVIT 343; and one of the WVIT 444 or WVIT 450; Senior standing; and WVIT major."
